$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip "Background Processing" flag on for the remaining data row.
$ws.Range("G2").Value = $true

# The second sample row (row 3) was a duplicate kept only to exercise the
# "Background Processing = TRUE" test plan; consolidate down to one row now
# that G2 itself carries TRUE, removing the now-redundant row entirely.
$ws.Rows(3).Delete()

# Leave the whole (now-last) data row selected, as it was after the edit.
$ws.Rows(2).Select()
